$wb = $excel.ActiveWorkbook

# Source sheet to copy the layout/format from (last existing sheet, 20201013).
$src = $wb.Worksheets.Item("20201013")

# Duplicate it right after itself -- this is how the new sheet ended up with
# the same column widths / alternating style artifacts as the original.
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item($src.Index + 1)
$newSheet.Name = "20201020"

# The source sheet had column C manually widened to fit its long SQL
# formulas; the new sheet never got that treatment, so put column C back to
# the sheet's normal default width (only column B keeps a custom width here).
$defaultWidth = $src.Columns("D").ColumnWidth
$newSheet.Columns("C").ColumnWidth = $defaultWidth

# Rows 2-11 keep the same restaurant (da04f5c9...), ids bump from 139..148
# to 149..158 (the copied formulas in C2:C11 recompute automatically).
$newSheet.Range("A2").Value = 149
$newSheet.Range("A3").Value = 150
$newSheet.Range("A4").Value = 151
$newSheet.Range("A5").Value = 152
$newSheet.Range("A6").Value = 153
$newSheet.Range("A7").Value = 154
$newSheet.Range("A8").Value = 155
$newSheet.Range("A9").Value = 156
$newSheet.Range("A10").Value = 157
$newSheet.Range("A11").Value = 158

# Formatting for the newly-added rows 12-19 continues the same alternating
# style pattern already present in rows 2-11 (even row in col A = styled,
# odd row in col A = default; col B always styled).
$styledTemplate = $newSheet.Range("A2")
$plainTemplate = $newSheet.Range("A3")
$colBTemplate = $newSheet.Range("B2")

function Copy-RowStyle($ws, $row) {
    $colBTemplate.Copy() | Out-Null
    $ws.Range("B$row").PasteSpecial(-4122) | Out-Null
    if ($row % 2 -eq 0) {
        $styledTemplate.Copy() | Out-Null
    } else {
        $plainTemplate.Copy() | Out-Null
    }
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0
}

function Set-SqlFormula($ws, $row) {
    $formula = '=_xlfn.CONCAT("INSERT INTO photos(restaurant_id, name, type) VALUES(UuidToBin(''",B' + $row + ',"''), LPAD(",A' + $row + ',", 7, ''0''), ''dish''",");")'
    $ws.Range("C$row").Formula = $formula
}

# New row 12: same restaurant, id 159.
Copy-RowStyle $newSheet 12
$newSheet.Range("A12").Value = 159
$newSheet.Range("B12").Value = "da04f5c9-ffb0-11ea-ba65-065a10bcba76"
Set-SqlFormula $newSheet 12

# Rows 13-17: new restaurant da053615-ffb0-11ea-ba65-065a10bcba76, ids 1..5.
Copy-RowStyle $newSheet 13
$newSheet.Range("A13").Value = 1
$newSheet.Range("B13").Value = "da053615-ffb0-11ea-ba65-065a10bcba76"
Set-SqlFormula $newSheet 13

Copy-RowStyle $newSheet 14
$newSheet.Range("A14").Value = 2
$newSheet.Range("B14").Value = "da053615-ffb0-11ea-ba65-065a10bcba76"
Set-SqlFormula $newSheet 14

Copy-RowStyle $newSheet 15
$newSheet.Range("A15").Value = 3
$newSheet.Range("B15").Value = "da053615-ffb0-11ea-ba65-065a10bcba76"
Set-SqlFormula $newSheet 15

Copy-RowStyle $newSheet 16
$newSheet.Range("A16").Value = 4
$newSheet.Range("B16").Value = "da053615-ffb0-11ea-ba65-065a10bcba76"
Set-SqlFormula $newSheet 16

Copy-RowStyle $newSheet 17
$newSheet.Range("A17").Value = 5
$newSheet.Range("B17").Value = "da053615-ffb0-11ea-ba65-065a10bcba76"
Set-SqlFormula $newSheet 17

# Rows 18-19: restaurant 146da67c-0526-11eb-ba65-065a10bcba76, ids 5..6.
Copy-RowStyle $newSheet 18
$newSheet.Range("A18").Value = 5
$newSheet.Range("B18").Value = "146da67c-0526-11eb-ba65-065a10bcba76"
Set-SqlFormula $newSheet 18

Copy-RowStyle $newSheet 19
$newSheet.Range("A19").Value = 6
$newSheet.Range("B19").Value = "146da67c-0526-11eb-ba65-065a10bcba76"
Set-SqlFormula $newSheet 19

# Fix up the view state: old active sheet is no longer selected/tab-active,
# new sheet is now the tab-selected / active one with selection B9.
$src.Range("C2").Select()
$newSheet.Activate()
$newSheet.Range("B9").Select()
